$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2023-04-23 Sunday"

# Update the multiplication problems in the table (20 rows x 5 cols)
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "77×30="
$t.Cell(1,2).Range.Text = "93×95="
$t.Cell(1,3).Range.Text = "73×35="
$t.Cell(1,4).Range.Text = "65×60="
$t.Cell(1,5).Range.Text = "31×19="
$t.Cell(2,1).Range.Text = "15×46="
$t.Cell(2,2).Range.Text = "70×62="
$t.Cell(2,3).Range.Text = "18×40="
$t.Cell(2,4).Range.Text = "49×15="
$t.Cell(2,5).Range.Text = "45×27="
$t.Cell(3,1).Range.Text = "41×39="
$t.Cell(3,2).Range.Text = "98×66="
$t.Cell(3,3).Range.Text = "24×38="
$t.Cell(3,4).Range.Text = "97×16="
$t.Cell(3,5).Range.Text = "92×36="
$t.Cell(4,1).Range.Text = "42×92="
$t.Cell(4,2).Range.Text = "31×32="
$t.Cell(4,3).Range.Text = "99×23="
$t.Cell(4,4).Range.Text = "52×46="
$t.Cell(4,5).Range.Text = "23×75="
$t.Cell(5,1).Range.Text = "53×80="
$t.Cell(5,2).Range.Text = "24×41="
$t.Cell(5,3).Range.Text = "26×77="
$t.Cell(5,4).Range.Text = "13×43="
$t.Cell(5,5).Range.Text = "77×39="
$t.Cell(6,1).Range.Text = "81×90="
$t.Cell(6,2).Range.Text = "49×65="
$t.Cell(6,3).Range.Text = "86×90="
$t.Cell(6,4).Range.Text = "33×66="
$t.Cell(6,5).Range.Text = "25×26="
$t.Cell(7,1).Range.Text = "32×52="
$t.Cell(7,2).Range.Text = "68×72="
$t.Cell(7,3).Range.Text = "91×74="
$t.Cell(7,4).Range.Text = "65×46="
$t.Cell(7,5).Range.Text = "33×10="
$t.Cell(8,1).Range.Text = "50×74="
$t.Cell(8,2).Range.Text = "79×46="
$t.Cell(8,3).Range.Text = "55×41="
$t.Cell(8,4).Range.Text = "67×95="
$t.Cell(8,5).Range.Text = "74×85="
$t.Cell(9,1).Range.Text = "98×49="
$t.Cell(9,2).Range.Text = "40×82="
$t.Cell(9,3).Range.Text = "92×83="
$t.Cell(9,4).Range.Text = "23×33="
$t.Cell(9,5).Range.Text = "65×76="
$t.Cell(10,1).Range.Text = "42×63="
$t.Cell(10,2).Range.Text = "36×68="
$t.Cell(10,3).Range.Text = "84×77="
$t.Cell(10,4).Range.Text = "86×100="
$t.Cell(10,5).Range.Text = "20×97="
$t.Cell(11,1).Range.Text = "18×97="
$t.Cell(11,2).Range.Text = "27×66="
$t.Cell(11,3).Range.Text = "52×53="
$t.Cell(11,4).Range.Text = "31×69="
$t.Cell(11,5).Range.Text = "71×19="
$t.Cell(12,1).Range.Text = "91×89="
$t.Cell(12,2).Range.Text = "10×64="
$t.Cell(12,3).Range.Text = "63×45="
$t.Cell(12,4).Range.Text = "32×51="
$t.Cell(12,5).Range.Text = "85×37="
$t.Cell(13,1).Range.Text = "63×83="
$t.Cell(13,2).Range.Text = "56×11="
$t.Cell(13,3).Range.Text = "37×57="
$t.Cell(13,4).Range.Text = "21×95="
$t.Cell(13,5).Range.Text = "10×70="
$t.Cell(14,1).Range.Text = "47×62="
$t.Cell(14,2).Range.Text = "95×93="
$t.Cell(14,3).Range.Text = "13×85="
$t.Cell(14,4).Range.Text = "90×30="
$t.Cell(14,5).Range.Text = "88×97="
$t.Cell(15,1).Range.Text = "46×67="
$t.Cell(15,2).Range.Text = "66×47="
$t.Cell(15,3).Range.Text = "25×18="
$t.Cell(15,4).Range.Text = "80×14="
$t.Cell(15,5).Range.Text = "34×48="
$t.Cell(16,1).Range.Text = "17×38="
$t.Cell(16,2).Range.Text = "38×41="
$t.Cell(16,3).Range.Text = "57×26="
$t.Cell(16,4).Range.Text = "68×81="
$t.Cell(16,5).Range.Text = "25×41="
$t.Cell(17,1).Range.Text = "93×100="
$t.Cell(17,2).Range.Text = "64×19="
$t.Cell(17,3).Range.Text = "55×57="
$t.Cell(17,4).Range.Text = "37×90="
$t.Cell(17,5).Range.Text = "44×19="
$t.Cell(18,1).Range.Text = "74×53="
$t.Cell(18,2).Range.Text = "99×38="
$t.Cell(18,3).Range.Text = "11×75="
$t.Cell(18,4).Range.Text = "38×29="
$t.Cell(18,5).Range.Text = "87×33="
$t.Cell(19,1).Range.Text = "87×74="
$t.Cell(19,2).Range.Text = "82×17="
$t.Cell(19,3).Range.Text = "13×23="
$t.Cell(19,4).Range.Text = "85×46="
$t.Cell(19,5).Range.Text = "43×14="
$t.Cell(20,1).Range.Text = "68×66="
$t.Cell(20,2).Range.Text = "40×35="
$t.Cell(20,3).Range.Text = "47×53="
$t.Cell(20,4).Range.Text = "41×74="
$t.Cell(20,5).Range.Text = "93×100="
